$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 129
$ws.Range("I2").Value = 275
$ws.Range("J2").Value = 1322
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 403
$ws.Range("M2").Value = 24
$ws.Range("N2").Value = 227
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 160
$ws.Range("T2").Value = 238
$ws.Range("U2").Value = 24
$ws.Range("V2").Value = 2135
$ws.Range("X2").Value = 2124
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 28
$ws.Range("AA2").Value = 14
